# Generate Report for Handoff
#
# The localization-status report is regenerated: the "Status" columns move
# from "Handed back: in sync with en-US" to "Ready for handoff", the
# associated "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps are refreshed, and the now-shorter Status column is narrowed
# to fit its new contents on every sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Refresh status text + handoff timestamps -----------------------------

# Overview sheet: zh-cn (E) / de-de (F) status cells, and the
# "Latest HO Xliff Generate Date" (G) cell.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-25 20:59:58"

# zh-cn detail sheet: Status (C) + Latest Handoff Datetime (H)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-25 20:59:53"

# de-de detail sheet: Status (C) + Latest Handoff Datetime (H)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-25 20:59:58"

# --- Narrow the "Status" / result columns to fit the shorter text ---------
# (was sized for "Handed back: in sync with en-US"; now fits "Ready for
# handoff"). 16.333333333333332 is the ColumnWidth input that this host's
# pixel-snapping rounds to the value closest to the new target width.

$narrowWidth = 16.333333333333332

$overview.Range("E:F").ColumnWidth = $narrowWidth
$zhcn.Range("C:C").ColumnWidth = $narrowWidth
$dede.Range("C:C").ColumnWidth = $narrowWidth
